$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Inheritance: Classes inherit state and behavior (methods from their
#    superclasses)" paragraph - collapse the split runs (caused by the
#    spell-check proofErr wrapper around "superclasses") back into a single
#    plain run by rewriting the paragraph's text.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(4)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = ""
$ins = $d.Range($p.Range.Start, $p.Range.Start)
$ins.InsertAfter("Inheritance: Classes inherit state and behavior (methods from their superclasses)")

# ---------------------------------------------------------------------------
# 2) "Interface: Contract between a class and the outside workd. ..."
#    paragraph - same treatment, removing the proofErr wrapper around "workd".
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(5)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = ""
$ins = $d.Range($p.Range.Start, $p.Range.Start)
$ins.InsertAfter("Interface: Contract between a class and the outside workd. Class implements an interface and promises to provide the behavior published by that interface. ")

# ---------------------------------------------------------------------------
# 3) "Abstract classes: ..." paragraph gains a new trailing sentence, and the
#    "_GoBack" bookmark (previously at the end of the last paragraph) is
#    relocated to the end of this paragraph's text.
# ---------------------------------------------------------------------------

# Remove the old bookmark first (it will be recreated below at the new spot).
$bOld = $d.Bookmarks.Item("_GoBack")
$bOld.Delete()

$p = $d.Paragraphs.Item(6)
$insertPos = $p.Range.End - 1
$addition = ". Abstract methods are declared but does not contain any implementation. By using abstract classes, you can inherit the implementation of non-abstract classes. Interfaces all need implementation. "
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter($addition)

# Re-create "_GoBack" as a zero-length bookmark right after the text we just
# added (i.e. immediately before this paragraph's mark). Adding a bookmark
# directly at a zero-width "paragraph end - 1" position is unreliable, so a
# one-character placeholder is used as a safe anchor and then removed via the
# bookmark's own Range (collapsing it to zero width in place).
$afterPos = $insertPos + $addition.Length
$placeholder = $d.Range($afterPos, $afterPos)
$placeholder.InsertAfter("X")
$anchor = $d.Range($afterPos, $afterPos + 1)
$d.Bookmarks.Add("_GoBack", $anchor)
$bNew = $d.Bookmarks.Item("_GoBack")
$bNewRange = $bNew.Range
$bNewRange.Text = ""

# ---------------------------------------------------------------------------
# 4) Insert the new "State: ..." paragraph in place of the blank paragraph
#    that used to separate "Inheritance: Promotes..." from "Constructor:".
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(11)
$ins = $d.Range($p.Range.Start, $p.Range.Start)
$ins.InsertAfter("State: The set of values of the attributes of a particular object is called a state. ")
